# Fruta / hortaliza, semanal
# Insert a new weekly record row before row 33, shifting existing rows 33-46 down to 34-47.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 33 (pushes rows 33..46 down to 34..47)
$ws.Rows.Item(33).Insert()

# Populate the newly inserted row 33 with the new record values
$ws.Cells.Item(33, 1).Value = 10
$ws.Cells.Item(33, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(33, 3).Value = "La Araucanía"
$ws.Cells.Item(33, 4).Value = 44508
$ws.Cells.Item(33, 5).Value = 9
$ws.Cells.Item(33, 6).Value = 100112026
$ws.Cells.Item(33, 7).Value = "Haba"
$ws.Cells.Item(33, 8).Value = "Sin especificar"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 140
$ws.Cells.Item(33, 11).Value = 8000
$ws.Cells.Item(33, 12).Value = 8000
$ws.Cells.Item(33, 13).Value = 8000
$ws.Cells.Item(33, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(33, 15).Value = "Región del Maule"
$ws.Cells.Item(33, 16).Value = 320
$ws.Cells.Item(33, 17).Value = 25
$ws.Cells.Item(33, 18).Value = "Hortaliza"

# Ensure the date cell keeps the same date number format as the rest of the D column
$ws.Cells.Item(33, 4).NumberFormat = $ws.Cells.Item(34, 4).NumberFormat
